$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "868×7=6076" "121×6=726"
Replace-Text "711×5=3555" "123×6=738"
Replace-Text "205×8=1640" "857×9=7713"
Replace-Text "646×9=5814" "384×5=1920"
Replace-Text "678×3=2034" "273×5=1365"
Replace-Text "493×6=2958" "877×5=4385"
Replace-Text "723×8=5784" "226×5=1130"
Replace-Text "347×4=1388" "702×4=2808"
Replace-Text "596×7=4172" "252×2=504"
Replace-Text "171×8=1368" "680×3=2040"
Replace-Text "697×5=3485" "473×8=3784"
Replace-Text "204×9=1836" "239×9=2151"
Replace-Text "972×8=7776" "129×2=258"
Replace-Text "599×6=3594" "726×9=6534"
Replace-Text "235×4=940" "787×2=1574"
Replace-Text "466×6=2796" "826×7=5782"
Replace-Text "602×7=4214" "797×6=4782"
Replace-Text "828×8=6624" "274×2=548"
Replace-Text "609×3=1827" "803×5=4015"
Replace-Text "574×2=1148" "953×4=3812"
Replace-Text "181×6=1086" "291×3=873"
Replace-Text "601×6=3606" "803×3=2409"
Replace-Text "911×4=3644" "126×9=1134"
Replace-Text "475×4=1900" "176×9=1584"
Replace-Text "581×5=2905" "319×8=2552"
